$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Fix the stray double-space typo in the "Alicia Salvador" note (row 3,
#    column E) - "Supervised together with  Alicia Salvador" ->
#    "Supervised together with Alicia Salvador".
# ---------------------------------------------------------------------------
$e3 = $ws.Range("E3").Value2
$ws.Range("E3").Value = $e3.Replace("together with  Alicia", "together with Alicia")

# ---------------------------------------------------------------------------
# 2) Reorder the supervision entries by year: swap the two-row block for the
#    2013-2014 Stirling MSc (rows 6-7) with the two-row block for the
#    2019-2020 Bosque MSc (rows 8-9), so the 2019-2020 entry now comes first.
# ---------------------------------------------------------------------------
$cols = @("A", "B", "C", "D", "E")

# Snapshot the current contents of the two blocks before overwriting anything.
$row6 = @{}
$row7 = @{}
$row8 = @{}
$row9 = @{}
foreach ($col in $cols) {
    $row6[$col] = $ws.Range($col + "6").Value2
    $row7[$col] = $ws.Range($col + "7").Value2
    $row8[$col] = $ws.Range($col + "8").Value2
    $row9[$col] = $ws.Range($col + "9").Value2
}

foreach ($col in $cols) {
    $target6 = $row8[$col]
    if ($target6 -eq $null) { $ws.Range($col + "6").ClearContents() | Out-Null } else { $ws.Range($col + "6").Value = $target6 }

    $target7 = $row9[$col]
    if ($target7 -eq $null) { $ws.Range($col + "7").ClearContents() | Out-Null } else { $ws.Range($col + "7").Value = $target7 }

    $target8 = $row6[$col]
    if ($target8 -eq $null) { $ws.Range($col + "8").ClearContents() | Out-Null } else { $ws.Range($col + "8").Value = $target8 }

    $target9 = $row7[$col]
    if ($target9 -eq $null) { $ws.Range($col + "9").ClearContents() | Out-Null } else { $ws.Range($col + "9").Value = $target9 }
}

# Row heights travel with the moved blocks too: the manually-set 31.5pt
# custom height follows the Bosque entry to row 6, the old row 8/9 wrapped
# heights follow the Stirling entry down to rows 8/9, and row 7/9 swap their
# plain vs. wrapped heights accordingly.
$ws.Rows.Item(6).RowHeight = 31.5
$ws.Rows.Item(7).RowHeight = 30
$ws.Rows.Item(8).RowHeight = 30
$ws.Rows.Item(9).RowHeight = 15
$ws.Rows.Item(9).AutoFit() | Out-Null

# ---------------------------------------------------------------------------
# 3) Restore the active selection left behind by the editor.
# ---------------------------------------------------------------------------
$ws.Range("E14").Select() | Out-Null
